{"js": "// Apply yellow highlighting to the \"13\u00ba slide \u2013 Recursos Humanos \u2013\n// Pesquisadores \u2013 gr\u00e1fico 20 da pg 60\" paragraph, matching the highlight\n// already used on the neighbouring \"12\u00ba slide\" / \"gr\u00e1fico 18\" paragraph.\n\nconst body = context.document.body;\n\n// Locate the target paragraph via a distinctive, stable text fragment\n// (the page number after \"pg\" is the part most likely to vary/typo,\n// so we don't key the search on it).\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\"Recursos Humanos\") !== -1 && text.indexOf(\"gr\u00e1fico 20 da\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  // Search the body for the paragraph's exact text to get a Range that\n  // covers only the run content (no trailing paragraph mark) \u2014 setting\n  // font.highlightColor on this range turns every run inside it yellow\n  // without touching the paragraph-mark run properties.\n  const results = body.search(target.text, { matchCase: false });\n  results.load(\"items/text\");\n  await context.sync();\n\n  if (results.items.length > 0) {\n    results.items[0].font.highlightColor = \"Yellow\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply yellow highlighting to the \"13\u00ba slide \u2013 Recursos Humanos \u2013\n# Pesquisadores \u2013 gr\u00e1fico 20 da pg 60\" paragraph, matching the highlight\n# already used on the neighbouring \"12\u00ba slide\" / \"gr\u00e1fico 18\" paragraph.\n\n$d = $word.ActiveDocument\n\n# Locate the target paragraph via a distinctive, stable text fragment\n# (the page number after \"pg\" is the part most likely to vary/typo, so\n# we don't key the search on it).\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n    if ($t -like \"*Recursos Humanos*\" -and $t -like \"*gr\u00e1fico 20 da*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($target -ne $null) {\n    # Highlight the paragraph's text range (all of its runs) in yellow,\n    # without touching the paragraph mark itself.\n    $target.Range.HighlightColorIndex = \"wdYellow\"\n}\n\nWrite-Output \"done\"\n"}
